$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D/E price & volume cells to remain Text (matches original
# inlineStr cell type) even when the new value looks numeric, by using the
# leading apostrophe text-prefix that Excel's COM Range.Value setter honors.

$ws.Range("D2").Value = "'58.814.90"
$ws.Range("E2").Value = "  -0.84%  "
$ws.Range("D3").Value = "'2.493.45"
$ws.Range("E3").Value = "  -0.92%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'536.16"
$ws.Range("E5").Value = "  -0.18%  "
$ws.Range("D6").Value = "'136.49"
$ws.Range("E6").Value = "  -2.10%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.566"
$ws.Range("E8").Value = "  +0.34%  "
$ws.Range("D9").Value = "'2.515.10"
$ws.Range("E9").Value = "  -0.19%  "
$ws.Range("E10").Value = "  -0.22%  "
$ws.Range("E11").Value = "  -2.65%  "
$ws.Range("E12").Value = "  -1.70%  "
$ws.Range("E13").Value = "  -3.10%  "
$ws.Range("D14").Value = "'2.940.50"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "'22.92"
$ws.Range("E15").Value = "  -2.37%  "
$ws.Range("D16").Value = "'58.740.61"
$ws.Range("E16").Value = "  -0.81%  "
$ws.Range("E17").Value = "  -1.64%  "
$ws.Range("D18").Value = "'2.510.29"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("D19").Value = "'11.09"
$ws.Range("E19").Value = "  -0.56%  "
$ws.Range("D20").Value = "'4.27"
$ws.Range("E20").Value = "  -0.91%  "
$ws.Range("D21").Value = "'323.05"
$ws.Range("E21").Value = "  -0.54%  "
$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  -0.05%  "
$ws.Range("D23").Value = "'5.91"
$ws.Range("E23").Value = "  +1.57%  "
$ws.Range("D24").Value = "'65.21"
$ws.Range("E24").Value = "  +2.75%  "
$ws.Range("D25").Value = "'0.419"
$ws.Range("E25").Value = "  -1.99%  "
$ws.Range("D26").Value = "'0.166"
$ws.Range("E26").Value = "  -1.00%  "
$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  -2.95%  "
$ws.Range("D29").Value = "'6.69"
$ws.Range("E29").Value = "  -3.84%  "
$ws.Range("D30").Value = "'0.0₃0765"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("E31").Value = "  -1.22%  "
$ws.Range("D32").Value = "'166.25"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("E33").Value = "  +5.06%  "
$ws.Range("E34").Value = "  +1.75%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'18.38"
$ws.Range("E36").Value = "  -0.74%  "
$ws.Range("D37").Value = "'4.09"
$ws.Range("E37").Value = "  -4.25%  "
$ws.Range("E38").Value = "  -3.13%  "
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("D40").Value = "'0.811"
$ws.Range("E40").Value = "  -0.55%  "
$ws.Range("D41").Value = "'3.61"
$ws.Range("E41").Value = "  -2.12%  "
$ws.Range("D42").Value = "'284.95"
$ws.Range("E42").Value = "  +1.86%  "
$ws.Range("D43").Value = "'5.18"
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("D44").Value = "'132.83"
$ws.Range("E44").Value = "  +8.10%  "
$ws.Range("D45").Value = "'0.995"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("D46").Value = "'0.602"
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("D47").Value = "'10.87"
$ws.Range("E47").Value = "  +0.17%  "
$ws.Range("D48").Value = "'0.0927"
$ws.Range("E48").Value = "  -0.31%  "
$ws.Range("D49").Value = "'0.0506"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("D51").Value = "'17.22"
$ws.Range("E51").Value = "  -3.49%  "
